# Generate Report for Handback
# Update the Correspond Handoff / Handback datetimes on the per-language
# report sheets to reflect the newly generated handback timestamps.

$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D5").Value = "2016-02-18 04:07:55"
$wsZh.Range("G5").Value = "2016-02-18 04:08:39"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D5").Value = "2016-02-18 04:08:07"
$wsDe.Range("G5").Value = "2016-02-18 04:09:01"
